$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Label" header in column H, matching the style of the other headers
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# Fill in the Label column: 0 for Control rows, 1 for MDD rows
$labels = @(0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
